$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.069.49"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.404.82"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.22"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.62"
$ws.Range("E6").Value = "  +2.54%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.401.57"
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.549"
$ws.Range("E9").Value = "  -5.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.30"
$ws.Range("E10").Value = "  +1.35%  "
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.421"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.988.81"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.87"
$ws.Range("E15").Value = "  -2.12%  "
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.073.25"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.393.92"
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.11"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.43"
$ws.Range("E20").Value = "  -2.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "374.12"
$ws.Range("E21").Value = "  -1.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.79"
$ws.Range("E22").Value = "  -1.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.42"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.511"
$ws.Range("E25").Value = "  -3.24%  "
$ws.Range("E26").Value = "  -4.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.49"
$ws.Range("E27").Value = "  -4.39%  "
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.05"
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("E31").Value = "  -2.91%  "
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.77"
$ws.Range("E34").Value = "  -1.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.03"
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("E36").Value = "  -6.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.33"
$ws.Range("E37").Value = "  -1.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.854"
$ws.Range("E38").Value = "  +7.56%  "
$ws.Range("E39").Value = "  -2.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0725"
$ws.Range("E40").Value = "  -2.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "25.77"
$ws.Range("E41").Value = "  -1.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.73"
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.723.38"
$ws.Range("E43").Value = "  -5.44%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.41"
$ws.Range("E44").Value = "  -0.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.71"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.36"
$ws.Range("E46").Value = "  -3.59%  "
$ws.Range("E47").Value = "  -1.68%  "
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "329.52"
$ws.Range("E49").Value = "  +2.78%  "
$ws.Range("E50").Value = "  -2.50%  "
$ws.Range("E51").Value = "  -1.92%  "
